# Update "想去人数" (F) and "最低票价" (G) figures on both the "展览"
# and "全部类型" sheets (they carry duplicate data in this workbook).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("G2").Value = 68

    $ws.Range("F3").Value = 613
    $ws.Range("F4").Value = 2172
    $ws.Range("F6").Value = 12716
    $ws.Range("F10").Value = 465
    $ws.Range("F11").Value = 1157
    $ws.Range("F12").Value = 958
    $ws.Range("F13").Value = 13678
    $ws.Range("F14").Value = 14061
    $ws.Range("F19").Value = 19
    $ws.Range("F27").Value = 5176
    $ws.Range("F29").Value = 262
}
